# Updated symbol list on Wed Feb  8 10:18:29 UTC 2023 with GitHub Actions
# Applies refreshed price/volume/hour data to the crypto listing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="330.47"; E="-0.04%"; G="10"},
    @{Row=3; D="45.35"; E="2.26%"; G="10"},
    @{Row=4; D="5.596"; E="2.58%"; G="10"},
    @{Row=5; D="0.08343"; E="3.90%"; G="10"},
    @{Row=6; D="2.087"; E="4.80%"; G="10"},
    @{Row=7; D="0.9721"; E="1.94%"; G="10"},
    @{Row=8; D="2.544"; E="-0.79%"; G="10"},
    @{Row=9; D="0.1178"; E="3.43%"; G="10"},
    @{Row=10; D="0.1924"; E="1.39%"; G="10"},
    @{Row=11; D="10.30"; E="-3.25%"; G="10"},
    @{Row=12; D="0.09843"; E="-1.18%"; G="10"},
    @{Row=13; D="0.04680"; E="-2.84%"; G="10"},
    @{Row=14; D="0.1061"; E="-0.24%"; G="10"},
    @{Row=15; D="0.001291"; E="1.81%"; G="10"},
    @{Row=16; D="0.006074"; E="2.17%"; G="10"},
    @{Row=17; D="3.375"; E="0.27%"; G="10"},
    @{Row=18; D="4.452"; E="1.26%"; G="10"},
    @{Row=19; D="0.3343"; E="-2.48%"; G="10"},
    @{Row=20; D=$null; E="-0.93%"; G="10"},
    @{Row=21; D="0.2639"; E="5.50%"; G="10"},
    @{Row=22; D="0.04179"; E="2.61%"; G="10"},
    @{Row=23; D="0.001302"; E="2.36%"; G="10"},
    @{Row=24; D="0.004552"; E="4.50%"; G="10"},
    @{Row=25; D=$null; E="8.78%"; G="10"},
    @{Row=26; D="0.0003751"; E="0.35%"; G="10"},
    @{Row=27; D=$null; E=$null; G="10"},
    @{Row=28; D=$null; E=$null; G="10"},
    @{Row=29; D=$null; E=$null; G="10"},
    @{Row=30; D=$null; E=$null; G="10"},
    @{Row=31; D=$null; E=$null; G="10"},
    @{Row=32; D=$null; E=$null; G="10"},
    @{Row=33; D=$null; E=$null; G="10"},
    @{Row=34; D=$null; E=$null; G="10"},
    @{Row=35; D=$null; E=$null; G="10"},
    @{Row=36; D=$null; E=$null; G="10"},
    @{Row=37; D=$null; E=$null; G="10"},
    @{Row=38; D="0.02698"; E="3.71%"; G="10"},
    @{Row=39; D="0.05759"; E="-0.63%"; G="10"},
    @{Row=40; D="0.007850"; E="4.13%"; G="10"},
    @{Row=41; D="0.1432"; E="2.04%"; G="10"},
    @{Row=42; D="0.007306"; E="-0.49%"; G="10"},
    @{Row=43; D="0.002129"; E="5.82%"; G="10"},
    @{Row=44; D="0.008513"; E="-3.68%"; G="10"},
    @{Row=45; D="0.3536"; E=$null; G="10"},
    @{Row=46; D="0.00007126"; E="0.94%"; G="10"},
    @{Row=47; D=$null; E="0.42%"; G="10"},
    @{Row=48; D="0.0005818"; E="0.44%"; G="10"},
    @{Row=49; D="0.003490"; E="-0.11%"; G="10"},
    @{Row=50; D="0.003509"; E="-0.45%"; G="10"},
    @{Row=51; D=$null; E="0.42%"; G="10"}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
    if ($null -ne $u.G) {
        $cell = $ws.Cells.Item($r, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $u.G
    }
}
